# Fruta / hortaliza, semanal
# Insert a new weekly record as row 6, shifting the existing rows 6-22 down to 7-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6 (this also copies the cell formatting - e.g. the
# date number format used in column D - from the row above, which is what we want).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44614
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100104
$ws.Cells.Item(6, 8).Value = "Frutos de pepita"
$ws.Cells.Item(6, 9).Value = 100104001
$ws.Cells.Item(6, 10).Value = "Granada"
$ws.Cells.Item(6, 11).Value = "Wonderfull"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 54
$ws.Cells.Item(6, 14).Value = 14000
$ws.Cells.Item(6, 15).Value = 14000
$ws.Cells.Item(6, 16).Value = 14000
$ws.Cells.Item(6, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(6, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 19).Value = 1000
$ws.Cells.Item(6, 20).Value = 14
